# Refresh the NATMI ligand-receptor edge table (Ccl12 -> Ccr3) with newly
# recomputed TPM-based statistics.
#
# The new run introduces a 4th sending-cluster category ("ECs") in addition
# to the existing Inflammatory-Mac / Neutrophils / Resolving-Mac clusters, so
# every "Sending cluster" x "Target cluster" combination (4 x 3 = 12 rows) is
# rewritten: rows 2-10 get updated values and rows 11-13 are newly added.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Row => values for columns A.. T (Sending cluster, Ligand symbol, Receptor
# symbol, Target cluster, then the 16 numeric NATMI statistics columns).
$rowData = @{
    2 = @("ECs", "Ccl12", "Ccr3", "Inflammatory-Mac", 1, 0.3333333333333333, 1.005755333333333, 3.017266, 0.01048729000197281, 0.01048729000197281, 2, 0.6666666666666666, 0.1808983333333333, 0.5426949999999999, 0.09546831801815302, 0.09546831801815302, 0.1819394635411111, 1.63745517187, 0.001001203937056936, 0.001001203937056936)
    3 = @("ECs", "Ccl12", "Ccr3", "Neutrophils", 1, 0.3333333333333333, 1.005755333333333, 3.017266, 0.01048729000197281, 0.01048729000197281, 3, 1, 1.572737, 4.718211, 0.8300051930177132, 0.8300051930177132, 1.581788625680667, 14.236097631126, 0.008704505162320172, 0.008704505162320172)
    4 = @("ECs", "Ccl12", "Ccr3", "Resolving-Mac", 1, 0.3333333333333333, 1.005755333333333, 3.017266, 0.01048729000197281, 0.01048729000197281, 1, 0.3333333333333333, 0.1412166666666667, 0.42365, 0.0745264889641337, 0.07452648896413369, 0.1420294156555556, 1.2782647409, 0.000781580902595696, 0.0007815809025956959)
    5 = @("Inflammatory-Mac", "Ccl12", "Ccr3", "Inflammatory-Mac", 3, 1, 38.755371, 116.266113, 0.4041130097356814, 0.4041130097356814, 2, 0.6666666666666666, 0.1808983333333333, 0.5426949999999999, 0.09546831801815302, 0.09546831801815302, 7.010782021614999, 63.09703819453499, 0.038579989328719, 0.038579989328719)
    6 = @("Inflammatory-Mac", "Ccl12", "Ccr3", "Neutrophils", 3, 1, 38.755371, 116.266113, 0.4041130097356814, 0.4041130097356814, 3, 1, 1.572737, 4.718211, 0.8300051930177132, 0.8300051930177132, 60.952005920427, 548.5680532838429, 0.3354158966466332, 0.3354158966466332)
    7 = @("Inflammatory-Mac", "Ccl12", "Ccr3", "Resolving-Mac", 3, 1, 38.755371, 116.266113, 0.4041130097356814, 0.4041130097356814, 1, 0.3333333333333333, 0.1412166666666667, 0.42365, 0.0745264889641337, 0.07452648896413369, 5.47290430805, 49.25613877245, 0.03011712376032911, 0.03011712376032911)
    8 = @("Neutrophils", "Ccl12", "Ccr3", "Inflammatory-Mac", 3, 1, 7.988471999999999, 23.965416, 0.08329801469605898, 0.08329801469605898, 2, 0.6666666666666666, 0.1808983333333333, 0.5426949999999999, 0.09546831801815302, 0.09546831801815302, 1.445101270679999, 13.00591143612, 0.007952321357284142, 0.007952321357284142)
    9 = @("Neutrophils", "Ccl12", "Ccr3", "Neutrophils", 3, 1, 7.988471999999999, 23.965416, 0.08329801469605898, 0.08329801469605898, 3, 1, 1.572737, 4.718211, 0.8300051930177132, 0.8300051930177132, 12.563765487864, 113.073889390776, 0.06913778476579474, 0.06913778476579474)
    10 = @("Neutrophils", "Ccl12", "Ccr3", "Resolving-Mac", 3, 1, 7.988471999999999, 23.965416, 0.08329801469605898, 0.08329801469605898, 1, 0.3333333333333333, 0.1412166666666667, 0.42365, 0.0745264889641337, 0.07452648896413369, 1.1281053876, 10.1529484884, 0.006207908572980087, 0.006207908572980085)
    11 = @("Resolving-Mac", "Ccl12", "Ccr3", "Inflammatory-Mac", 3, 1, 48.15271133333334, 144.458134, 0.5021016855662869, 0.5021016855662868, 2, 0.6666666666666666, 0.1808983333333333, 0.5426949999999999, 0.09546831801815302, 0.09546831801815302, 8.71074522568111, 78.39670703112999, 0.04793480339509294, 0.04793480339509294)
    12 = @("Resolving-Mac", "Ccl12", "Ccr3", "Neutrophils", 3, 1, 48.15271133333334, 144.458134, 0.5021016855662869, 0.5021016855662868, 3, 1, 1.572737, 4.718211, 0.8300051930177132, 0.8300051930177132, 75.73155076425267, 681.5839568782741, 0.4167470064429651, 0.416747006442965)
    13 = @("Resolving-Mac", "Ccl12", "Ccr3", "Resolving-Mac", 3, 1, 48.15271133333334, 144.458134, 0.5021016855662869, 0.5021016855662868, 1, 0.3333333333333333, 0.1412166666666667, 0.42365, 0.0745264889641337, 0.07452648896413369, 6.799965385455557, 61.1996884691, 0.03741987572822881, 0.03741987572822879)
}

foreach ($r in ($rowData.Keys | Sort-Object)) {
    $values = $rowData[$r]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $r).Value = $values[$i]
    }
}
